$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week 2")

# Row 11 - id 5
$ws.Range("C11").Value = 0.4375
$ws.Range("D11").Value = 0.44791666666666669
$ws.Range("F11").Value = "MouseRect in draw aangemaakt"

# Row 12 - id 6
$ws.Range("C12").Value = 0.4513888888888889
$ws.Range("D12").Value = 0.47569444444444442
$ws.Range("F12").Value = "Eerste knop werkend gemaakt met muis"

# Update active selection to F12
[void]$ws.Range("F12").Select()
